# Fix preset sync and menu ordering
# Rows 19 and 20 on the "Menu Mock" sheet had their Option/Control/Values/
# Source/Tooltip data swapped (Haptic Feedback <-> Dynamic Intensity).
# Category (col A) and Default (col D) are identical on both rows, so only
# columns B, C, E, F, G need to be exchanged between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menu Mock")

$cols = @("B", "C", "E", "F", "G")

foreach ($col in $cols) {
    $row19Value = $ws.Range("$col`19").Value2
    $row20Value = $ws.Range("$col`20").Value2

    $ws.Range("$col`19").Value2 = $row20Value
    $ws.Range("$col`20").Value2 = $row19Value
}
